{"js": "// R10-gasser18.docx \u2014 \"Updated language for reading 10\"\n//\n// Two independent textual tweaks in the body:\n//  1. \"I believe what makes for good indicator ...\" -> \"... what makes for\n//     a good indicator ...\" (insert the missing article \"a \").\n//  2. The quality-assessment/quality-assurance paragraph is re-punctuated\n//     and gets a new transition sentence:\n//       \"... already developed by relying on ... testing phases, while\n//        quality assurance cares about ...\"\n//     becomes\n//       \"... already developed, by relying on ... testing phases. On the\n//        other hand, quality assurance cares about ...\"\n\nconst body = context.document.body;\n\n// --- Edit 1: insert \"a \" before \"good indicator is the one that comes after\"\nconst hit1 = body.search(\"good indicator is the one that comes after\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\n\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\"a \", \"Before\");\n}\n\n// --- Edit 2: re-punctuate + add \"On the other hand\" transition\nconst oldText2 =\n  \"developed by relying on the outcome of the code development, unit test, and integration testing phases, while quality assurance cares about \";\nconst newText2 =\n  \"developed, by relying on the outcome of the code development, unit test, and integration testing phases. On the other hand, quality assurance cares about \";\n\nconst hit2 = body.search(oldText2, { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\n\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(newText2, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# R10-gasser18.docx \u2014 \"Updated language for reading 10\"\n#\n# Two independent textual tweaks in the body:\n#  1. \"I believe what makes for good indicator ...\" -> \"... what makes for\n#     a good indicator ...\" (insert the missing article \"a \").\n#  2. The quality-assessment/quality-assurance paragraph is re-punctuated\n#     and gets a new transition sentence:\n#       \"... already developed by relying on ... testing phases, while\n#        quality assurance cares about ...\"\n#     becomes\n#       \"... already developed, by relying on ... testing phases. On the\n#        other hand, quality assurance cares about ...\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: insert \"a \" before \"good indicator is the one that comes after\"\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"good indicator is the one that comes after\")\nif ($found1) {\n    $rng1.Collapse(1)\n    $rng1.InsertBefore(\"a \")\n}\n\n# --- Edit 2: re-punctuate + add \"On the other hand\" transition\n$oldText2 = \"developed by relying on the outcome of the code development, unit test, and integration testing phases, while quality assurance cares about \"\n$newText2 = \"developed, by relying on the outcome of the code development, unit test, and integration testing phases. On the other hand, quality assurance cares about \"\n\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute($oldText2)\nif ($found2) {\n    $rng2.Text = $newText2\n}\n"}
